$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply text-valued cell updates while preserving the default (un-styled) cell format.
function Set-TextValue {
    param($Worksheet, $CellRef, $NewValue)
    $rng = $Worksheet.Range($CellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $NewValue
    $rng.Style = "Normal"
}

Set-TextValue $ws "D2" "277.88"
Set-TextValue $ws "E2" "1.62%"
Set-TextValue $ws "D3" "27.20"
Set-TextValue $ws "E3" "1.65%"
Set-TextValue $ws "D4" "4.877"
Set-TextValue $ws "E4" "-0.57%"
Set-TextValue $ws "D5" "0.06369"
Set-TextValue $ws "E5" "0.67%"
Set-TextValue $ws "D6" "6.978"
Set-TextValue $ws "E6" "0.99%"
Set-TextValue $ws "D7" "1.253"
Set-TextValue $ws "E7" "-9.96%"
Set-TextValue $ws "D8" "0.8809"
Set-TextValue $ws "E8" "-0.27%"
Set-TextValue $ws "D9" "0.1525"
Set-TextValue $ws "E9" "3.52%"
Set-TextValue $ws "D10" "0.05130"
Set-TextValue $ws "E10" "0.59%"
Set-TextValue $ws "E11" "1.73%"
Set-TextValue $ws "D12" "0.02965"
Set-TextValue $ws "E12" "-6.90%"
Set-TextValue $ws "D13" "0.09009"
Set-TextValue $ws "E13" "-0.47%"
Set-TextValue $ws "D14" "0.001568"
Set-TextValue $ws "E14" "0.51%"
Set-TextValue $ws "D15" "0.0006395"
Set-TextValue $ws "E15" "0.96%"
Set-TextValue $ws "D16" "0.005889"
Set-TextValue $ws "E16" "-2.48%"
Set-TextValue $ws "D17" "3.461"
Set-TextValue $ws "E17" "-0.29%"
Set-TextValue $ws "D18" "3.320"
Set-TextValue $ws "E18" "-0.86%"
Set-TextValue $ws "E20" "0.92%"
Set-TextValue $ws "D21" "0.1337"
Set-TextValue $ws "E21" "0.23%"
Set-TextValue $ws "D22" "3.900"
Set-TextValue $ws "E22" "-0.67%"
Set-TextValue $ws "D23" "0.04419"
Set-TextValue $ws "E23" "1.96%"
Set-TextValue $ws "D24" "0.001172"
Set-TextValue $ws "E24" "-0.40%"
Set-TextValue $ws "D25" "0.003880"
Set-TextValue $ws "E25" "6.43%"
Set-TextValue $ws "D26" "0.0001200"
Set-TextValue $ws "E26" "-0.07%"
Set-TextValue $ws "E27" "13.96%"
Set-TextValue $ws "E40" "2.03%"
Set-TextValue $ws "D41" "0.006805"
Set-TextValue $ws "E41" "2.99%"
Set-TextValue $ws "D42" "0.1178"
Set-TextValue $ws "E42" "0.96%"
Set-TextValue $ws "D43" "0.002020"
Set-TextValue $ws "E43" "-9.08%"
Set-TextValue $ws "D44" "0.01123"
Set-TextValue $ws "E44" "-10.72%"
Set-TextValue $ws "D45" "0.00005187"
Set-TextValue $ws "E45" "-2.73%"
Set-TextValue $ws "D46" "1.129"
Set-TextValue $ws "E46" "-52.08%"
Set-TextValue $ws "D47" "0.02024"
Set-TextValue $ws "E47" "-4.59%"
